$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1450.8572
$ws.Range("I32").Value = 1117.6
$ws.Range("J32").Value = 1636
$ws.Range("K32").Value = 1117.6
$ws.Range("L32").Value = 1636
$ws.Range("M32").Value = -791.5999999999999
$ws.Range("N32").Value = -2288
$ws.Range("H33").Value = 128.85715
$ws.Range("I33").Value = 125.333336
$ws.Range("K33").Value = 125.333336
$ws.Range("M33").Value = 103.666664
$ws.Range("H40").Value = 2199.9473
$ws.Range("I40").Value = 1600
$ws.Range("J40").Value = 2312.4375
$ws.Range("K40").Value = 1600
$ws.Range("L40").Value = 2312.4375
$ws.Range("M40").Value = -1425
$ws.Range("N40").Value = -2662.4375
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 16024.5
$ws.Range("J44").Value = 16024.5
$ws.Range("L44").Value = 16024.5
$ws.Range("N44").Value = -17000.5
$ws.Range("H55").Value = 19406.9
$ws.Range("J55").Value = 19406.9
$ws.Range("L55").Value = 19406.9
$ws.Range("N55").Value = -20036.9
$ws.Range("H63").Value = 1466.6666
$ws.Range("I63").Value = 1538.4615
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 1538.4615
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = -852.4614999999999
$ws.Range("N63").Value = -2372
$ws.Range("H66").Value = 1466.6666
$ws.Range("I66").Value = 1538.4615
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 7692.307499999999
$ws.Range("L66").Value = 5000
$ws.Range("M66").Value = -4260.307499999999
$ws.Range("N66").Value = -11864
$ws.Range("H74").Value = 875.4838999999999
$ws.Range("I74").Value = 837.04
$ws.Range("J74").Value = 1035.6666
$ws.Range("K74").Value = 837.04
$ws.Range("L74").Value = 1035.6666
$ws.Range("M74").Value = 36.96000000000004
$ws.Range("N74").Value = -2783.6666
$ws.Range("H77").Value = 875.4838999999999
$ws.Range("I77").Value = 837.04
$ws.Range("J77").Value = 1035.6666
$ws.Range("K77").Value = 4185.2
$ws.Range("L77").Value = 5178.333000000001
$ws.Range("M77").Value = 182.8000000000002
$ws.Range("N77").Value = -13914.333
$ws.Range("H80").Value = 21076
$ws.Range("J80").Value = 21076
$ws.Range("L80").Value = 21076
$ws.Range("N80").Value = -23072
$ws.Range("H83").Value = 21076
$ws.Range("J83").Value = 21076
$ws.Range("L83").Value = 63228
$ws.Range("N83").Value = -73212

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 71468
$ws.Range("I20").Value = 101296.664
$ws.Range("J20").Value = 26725
$ws.Range("K20").Value = 101296.664
$ws.Range("L20").Value = 26725
$ws.Range("M20").Value = -101049.664
$ws.Range("N20").Value = -27219
$ws.Range("H82").Value = 16264.8
$ws.Range("I82").Value = 12000
$ws.Range("J82").Value = 18092.572
$ws.Range("K82").Value = 12000
$ws.Range("L82").Value = 18092.572
$ws.Range("M82").Value = -11617
$ws.Range("N82").Value = -18858.572
$ws.Range("H85").Value = 16264.8
$ws.Range("I85").Value = 12000
$ws.Range("J85").Value = 18092.572
$ws.Range("K85").Value = 12000
$ws.Range("L85").Value = 18092.572
$ws.Range("M85").Value = -10674
$ws.Range("N85").Value = -20744.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 23700
$ws.Range("J68").Value = 23700
$ws.Range("L68").Value = 23700
$ws.Range("N68").Value = -25198
$ws.Range("H71").Value = 23700
$ws.Range("J71").Value = 23700
$ws.Range("L71").Value = 71100
$ws.Range("N71").Value = -78588
$ws.Range("H74").Value = 12825.6
$ws.Range("J74").Value = 12825.6
$ws.Range("L74").Value = 12825.6
$ws.Range("N74").Value = -14573.6
$ws.Range("H77").Value = 12825.6
$ws.Range("J77").Value = 12825.6
$ws.Range("L77").Value = 38476.8
$ws.Range("N77").Value = -47212.8
$ws.Range("H99").Value = 2739.1304
$ws.Range("I99").Value = 2245.4546
$ws.Range("J99").Value = 3191.6667
$ws.Range("K99").Value = 2245.4546
$ws.Range("L99").Value = 3191.6667
$ws.Range("M99").Value = -747.4546
$ws.Range("N99").Value = -6187.6667
$ws.Range("H126").Value = 2739.1304
$ws.Range("I126").Value = 2245.4546
$ws.Range("J126").Value = 3191.6667
$ws.Range("K126").Value = 6736.3638
$ws.Range("L126").Value = 9575.000100000001
$ws.Range("M126").Value = -4266.3638
$ws.Range("N126").Value = -14515.0001
$ws.Range("H134").Value = 2477.6924
$ws.Range("I134").Value = 1773.4828
$ws.Range("J134").Value = 4519.9
$ws.Range("K134").Value = 5320.4484
$ws.Range("L134").Value = 13559.7
$ws.Range("M134").Value = -2785.4484
$ws.Range("N134").Value = -18629.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 250.41667
$ws.Range("I7").Value = 176.42857
$ws.Range("J7").Value = 354
$ws.Range("K7").Value = 529.28571
$ws.Range("L7").Value = 1062
$ws.Range("M7").Value = -417.28571
$ws.Range("N7").Value = -1286
$ws.Range("H92").Value = 375
$ws.Range("I92").Value = 150
$ws.Range("J92").Value = 487.5
$ws.Range("K92").Value = 450
$ws.Range("L92").Value = 1462.5
$ws.Range("M92").Value = 798
$ws.Range("N92").Value = -3958.5
$ws.Range("H107").Value = 486932.88
$ws.Range("I107").Value = 959
$ws.Range("J107").Value = 1556075.4
$ws.Range("K107").Value = 2877
$ws.Range("L107").Value = 4668226.199999999
$ws.Range("M107").Value = -957
$ws.Range("N107").Value = -4672066.199999999
$ws.Range("H131").Value = 813.05
$ws.Range("J131").Value = 820.8674
$ws.Range("L131").Value = 2462.6022
$ws.Range("N131").Value = -12542.6022

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3429.1667
$ws.Range("I80").Value = 3377.2727
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 3377.2727
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -2379.2727
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 3429.1667
$ws.Range("I83").Value = 3377.2727
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 16886.3635
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -11894.3635
$ws.Range("N83").Value = -29984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1885.5714
$ws.Range("I68").Value = 1759.8
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 1759.8
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -1010.8
$ws.Range("N68").Value = -3698
$ws.Range("H71").Value = 1885.5714
$ws.Range("I71").Value = 1759.8
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 8799
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -5055
$ws.Range("N71").Value = -18488
$ws.Range("H136").Value = 2443.5757
$ws.Range("I136").Value = 1252.12
$ws.Range("J136").Value = 6166.875
$ws.Range("K136").Value = 3756.36
$ws.Range("L136").Value = 18500.625
$ws.Range("M136").Value = -1206.36
$ws.Range("N136").Value = -23600.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1938.1428
$ws.Range("I81").Value = 1683.7778
$ws.Range("J81").Value = 2396
$ws.Range("K81").Value = 3367.5556
$ws.Range("L81").Value = 4792
$ws.Range("M81").Value = -2306.5556
$ws.Range("N81").Value = -6914
$ws.Range("H84").Value = 1938.1428
$ws.Range("I84").Value = 1683.7778
$ws.Range("J84").Value = 2396
$ws.Range("K84").Value = 16837.778
$ws.Range("L84").Value = 23960
$ws.Range("M84").Value = -11533.778
$ws.Range("N84").Value = -34568
$ws.Range("H86").Value = 12500
$ws.Range("J86").Value = 12500
$ws.Range("L86").Value = 12500
$ws.Range("N86").Value = -14746
$ws.Range("H89").Value = 12500
$ws.Range("J89").Value = 12500
$ws.Range("L89").Value = 62500
$ws.Range("N89").Value = -73732
